# mereni_MD/Prumery/Sval1.xlsx - "oprava exel tabulek done"
# The raw measurements in columns A (Sklon mV) and E (Sklon mBar) were
# recorded with a misplaced decimal point; correct them by dividing by 10.
# The AVERAGE formulas in row 2 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 4-9
$ws.Range("A4").Value = 2.3786999999999998
$ws.Range("A5").Value = 2.4123999999999999
$ws.Range("A6").Value = 2.3832
$ws.Range("A7").Value = 2.3910999999999998
$ws.Range("A8").Value = 2.3889999999999998
$ws.Range("A9").Value = 2.3647

# Column E, rows 4-9
$ws.Range("E4").Value = 0.38508999999999999
$ws.Range("E5").Value = 0.39173000000000002
$ws.Range("E6").Value = 0.38497999999999999
$ws.Range("E7").Value = 0.38874999999999998
$ws.Range("E8").Value = 0.38978000000000002
$ws.Range("E9").Value = 0.38430999999999998

# Move/select the active cell from M10 to E10
$ws.Range("E10").Select()
